# Update countries & provincias Spain
#
# The source COVID-19 dashboard refreshed its "Ciudad" (province) table:
#   - the timestamp banner in A1 moved from 13:22 to 13:52
#   - several provinces' rows (Valladolid/Castilla-La Mancha/Malaga/Salamanca/
#     Sevilla/Segovia/Asturias/Gipuzkoa-Guipuzcoa and Caceres/Cantabria and
#     Castello-Castellon/Jaen) were re-labelled because the underlying data
#     pull reordered those entries
#   - most rows got refreshed totals (Casos totales/activos, Recuperados,
#     Muertes) for the 13:52 snapshot

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 13:52"

# Row 8 - Ciudad Real
$ws.Range("B8").Value = 5717
$ws.Range("C8").Value = 2998
$ws.Range("D8").Value = 9927
$ws.Range("E8").Value = 659

# Row 9 - Valencia/Valencia
$ws.Range("B9").Value = 4849

# Row 11 - Albacete
$ws.Range("B11").Value = 3575
$ws.Range("C11").Value = 2998
$ws.Range("D11").Value = 9927
$ws.Range("E11").Value = 338

# Row 13 - Alacant/Alicante
$ws.Range("B13").Value = 3378

# Row 15 - Toledo
$ws.Range("B15").Value = 3193
$ws.Range("C15").Value = 2998
$ws.Range("D15").Value = 9927
$ws.Range("E15").Value = 454

# Row 17 - now Valladolid (was Castilla-La Mancha)
$ws.Range("A17").Value = "Valladolid"
$ws.Range("B17").Value = 2836
$ws.Range("C17").Value = 943
$ws.Range("D17").Value = 1670
$ws.Range("E17").Value = 223

# Row 18 - now Castilla-La Mancha (was Valladolid)
$ws.Range("A18").Value = "Castilla-La Mancha"
$ws.Range("B18").Value = 2780
$ws.Range("C18").Value = 71
$ws.Range("D18").Value = 2446
$ws.Range("E18").Value = 263

# Row 19 - now Malaga (was Salamanca)
$ws.Range("A19").Value = "Malaga"
$ws.Range("B19").Value = 2321
$ws.Range("C19").Value = 664
$ws.Range("D19").Value = 1456
$ws.Range("E19").Value = 201

# Row 20 - now Salamanca (was Malaga)
$ws.Range("A20").Value = "Salamanca"
$ws.Range("B20").Value = 2291
$ws.Range("C20").Value = 686
$ws.Range("D20").Value = 1338
$ws.Range("E20").Value = 267

# Row 21 - now Sevilla (was Asturias)
$ws.Range("A21").Value = "Sevilla"
$ws.Range("B21").Value = 2159
$ws.Range("C21").Value = 327
$ws.Range("D21").Value = 1650
$ws.Range("E21").Value = 182

# Row 22 - now Segovia (was Gipuzkoa/Guipuzcoa)
$ws.Range("A22").Value = "Segovia"
$ws.Range("B22").Value = 2103
$ws.Range("C22").Value = 582
$ws.Range("D22").Value = 1363
$ws.Range("E22").Value = 158

# Row 23 - now Asturias (was Sevilla)
$ws.Range("A23").Value = "Asturias"
$ws.Range("B23").Value = 2096
$ws.Range("C23").Value = 487
$ws.Range("D23").Value = 1443
$ws.Range("E23").Value = 166

# Row 24 - now Gipuzkoa/Guipuzcoa (was Segovia)
$ws.Range("A24").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B24").Value = 2086
$ws.Range("C24").Value = 5193
$ws.Range("D24").Value = 5174
$ws.Range("E24").Value = 136

# Row 26 - Leon
$ws.Range("B26").Value = 1935
$ws.Range("C26").Value = 906
$ws.Range("D26").Value = 756
$ws.Range("E26").Value = 273

# Row 27 - Granada
$ws.Range("B27").Value = 1911
$ws.Range("C27").Value = 409
$ws.Range("D27").Value = 1321
$ws.Range("E27").Value = 181

# Row 28 - now Caceres (was Cantabria)
$ws.Range("A28").Value = "Caceres"
$ws.Range("B28").Value = 1860
$ws.Range("C28").Value = 290
$ws.Range("D28").Value = 1304
$ws.Range("E28").Value = 266

# Row 29 - now Cantabria (was Caceres)
$ws.Range("A29").Value = "Cantabria"
$ws.Range("B29").Value = 1823
$ws.Range("C29").Value = 363
$ws.Range("D29").Value = 1328
$ws.Range("E29").Value = 132

# Row 32 - Burgos
$ws.Range("B32").Value = 1320
$ws.Range("C32").Value = 564
$ws.Range("D32").Value = 607
$ws.Range("E32").Value = 149

# Row 34 - Guadalajara
$ws.Range("B34").Value = 1226
$ws.Range("C34").Value = 2998
$ws.Range("D34").Value = 9927
$ws.Range("E34").Value = 157

# Row 35 - Cordoba
$ws.Range("B35").Value = 1211
$ws.Range("C35").Value = 207
$ws.Range("D35").Value = 945
$ws.Range("E35").Value = 59

# Row 36 - now Castello/Castellon (was Jaen)
$ws.Range("A36").Value = "Castello/Castellon"
$ws.Range("B36").Value = 1192
$ws.Range("C36").Value = 274
$ws.Range("D36").Value = 736
$ws.Range("E36").Value = 123

# Row 37 - now Jaen (was Castello/Castellon)
$ws.Range("A37").Value = "Jaen"
$ws.Range("B37").Value = 1189
$ws.Range("C37").Value = 220
$ws.Range("D37").Value = 854
$ws.Range("E37").Value = 115

# Row 38 - Soria
$ws.Range("B38").Value = 1095
$ws.Range("C38").Value = 265
$ws.Range("D38").Value = 740
$ws.Range("E38").Value = 90

# Row 39 - Cadiz
$ws.Range("B39").Value = 1044
$ws.Range("C39").Value = 235
$ws.Range("D39").Value = 747
$ws.Range("E39").Value = 62

# Row 40 - Avila
$ws.Range("B40").Value = 987
$ws.Range("C40").Value = 402
$ws.Range("D40").Value = 481
$ws.Range("E40").Value = 104

# Row 41 - Cuenca
$ws.Range("B41").Value = 969
$ws.Range("C41").Value = 2998
$ws.Range("D41").Value = 9927
$ws.Range("E41").Value = 147

# Row 43 - Badajoz
$ws.Range("B43").Value = 902
$ws.Range("C43").Value = 325
$ws.Range("D43").Value = 501
$ws.Range("E43").Value = 76

# Row 45 - Palencia
$ws.Range("B45").Value = 648
$ws.Range("C45").Value = 195
$ws.Range("D45").Value = 400
$ws.Range("E45").Value = 53

# Row 49 - Zamora
$ws.Range("B49").Value = 482
$ws.Range("C49").Value = 182
$ws.Range("D49").Value = 245
$ws.Range("E49").Value = 55

# Row 51 - Almeria
$ws.Range("B51").Value = 416
$ws.Range("C51").Value = 100
$ws.Range("D51").Value = 278

# Row 52 - Huelva
$ws.Range("B52").Value = 344
$ws.Range("C52").Value = 62
$ws.Range("D52").Value = 255
